$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "36.430.46"
$ws.Range("E2").Value = "  -2.67%  "
Set-TextValue "D3" "1.981.46"
$ws.Range("E3").Value = "  -3.13%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "245.88"
$ws.Range("E5").Value = "  -2.56%  "
Set-TextValue "D6" "0.625"
$ws.Range("E6").Value = "  -4.05%  "
Set-TextValue "D7" "59.23"
$ws.Range("E7").Value = "  -9.55%  "
$ws.Range("E8").Value = "  +0.01%  "
Set-TextValue "D9" "0.375"
$ws.Range("E9").Value = "  -7.61%  "
Set-TextValue "D10" "56.84"
$ws.Range("E10").Value = "  -4.80%  "
Set-TextValue "D11" "0.0874"
$ws.Range("E11").Value = "  +10.35%  "
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D13" "0.859"
$ws.Range("E13").Value = "  -6.85%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D14" "22.44"
$ws.Range("E14").Value = "  -4.63%  "
Set-TextValue "D15" "2.269.38"
$ws.Range("E15").Value = "  -3.18%  "
Set-TextValue "D16" "13.83"
$ws.Range("E16").Value = "  -6.78%  "
Set-TextValue "D17" "5.47"
$ws.Range("E17").Value = "  -5.02%  "
Set-TextValue "D18" "1.973.86"
$ws.Range("E18").Value = "  -3.48%  "
Set-TextValue "D19" "36.311.88"
$ws.Range("E19").Value = "  -2.68%  "
Set-TextValue "D20" "0.0₃0907"
$ws.Range("E20").Value = "  +2.50%  "
Set-TextValue "D21" "70.49"
$ws.Range("E21").Value = "  -4.30%  "
$ws.Range("E22").Value = "  -4.19%  "
Set-TextValue "D23" "234.47"
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -4.13%  "
Set-TextValue "D26" "2.30"
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("E27").Value = "  -2.23%  "
Set-TextValue "D28" "165.28"
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D29" "19.94"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D30" "0.132"
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("E31").Value = "  -2.24%  "
Set-TextValue "D32" "1.19"
$ws.Range("E32").Value = "  -0.38%  "
Set-TextValue "D33" "4.89"
$ws.Range("E33").Value = "  -5.54%  "
Set-TextValue "D34" "0.0649"
$ws.Range("E34").Value = "  +3.14%  "
$ws.Range("E35").Value = "  -5.57%  "
$ws.Range("E36").Value = "  -0.01%  "
Set-TextValue "D37" "6.09"
$ws.Range("E37").Value = "  -4.09%  "
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("E39").Value = "  -6.86%  "
Set-TextValue "D40" "2.95"
$ws.Range("E40").Value = "  -2.47%  "
Set-TextValue "D41" "1.22"
$ws.Range("E41").Value = "  -5.45%  "
Set-TextValue "D42" "0.0966"
$ws.Range("E42").Value = "  -5.33%  "
Set-TextValue "D43" "2.89"
$ws.Range("E43").Value = "  -5.31%  "
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D45" "16.29"
$ws.Range("E45").Value = "  -6.89%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D46" "1.08"
$ws.Range("E46").Value = "  -6.50%  "
Set-TextValue "D47" "91.33"
$ws.Range("E47").Value = "  -4.93%  "
Set-TextValue "D48" "1.367.17"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("E49").Value = "  -5.18%  "
Set-TextValue "D50" "2.84"
$ws.Range("E50").Value = "  -2.66%  "
Set-TextValue "D51" "45.37"
$ws.Range("E51").Value = "  -4.30%  "

Write-Host "Applied 98 cell changes"
